$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 (Idaho) failed with a timeout on this run, so the scraped
# numeric results (Date Published .. Pct Deaths Black/AA) are blanked
# out instead of populated, just like the other error rows in the sheet.
# Writing a bare quote-prefix forces an empty *text* cell (rather than
# simply deleting the cell), then ClearFormats() drops the transient
# quote-prefix/number-format style so the cell ends up unstyled, same
# as the untouched error rows elsewhere in the sheet.
$ws.Range("B38:H38").Value = "'"
$ws.Range("B38:H38").ClearFormats()

# Pct Includes Unknown Race stays False; Pct Includes Hispanic Black
# flips from True to False now that there is no successful result.
$ws.Cells.Item(38, 9).Value = $false
$ws.Cells.Item(38, 10).Value = $false

# Status column records the failure instead of "Success!".
$ws.Cells.Item(38, 15).Value = "An error occurred. ... TimeoutException('', None, None)"
